$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 116; existing rows 116-184 shift down to 117-185.
$ws.Rows(116).Insert()

# Populate the newly inserted row 116 with the new record's data.
$ws.Cells.Item(116,1).Value = 5
$ws.Cells.Item(116,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(116,3).Value = "Maule"
$ws.Cells.Item(116,4).Value = 44606
$ws.Cells.Item(116,5).Value = 7
$ws.Cells.Item(116,6).Value = 100112021
$ws.Cells.Item(116,7).Value = "Ají"
$ws.Cells.Item(116,8).Value = "Americana (o)"
$ws.Cells.Item(116,9).Value = "Primera"
$ws.Cells.Item(116,10).Value = 150
$ws.Cells.Item(116,11).Value = 15000
$ws.Cells.Item(116,12).Value = 15000
$ws.Cells.Item(116,13).Value = 15000
$ws.Cells.Item(116,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(116,15).Value = "Región del Maule"
$ws.Cells.Item(116,16).Value = 600
$ws.Cells.Item(116,17).Value = 25
$ws.Cells.Item(116,18).Value = "Hortaliza"
